# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3863.037
$ws.Range("I51").Value = 4040.2
$ws.Range("K51").Value = 4040.2
$ws.Range("M51").Value = -3556.2
$ws.Range("H55").Value = 164.36363
$ws.Range("I55").Value = 101
$ws.Range("J55").Value = 333.33334
$ws.Range("K55").Value = 101
$ws.Range("L55").Value = 333.33334
$ws.Range("M55").Value = 113
$ws.Range("N55").Value = -761.33334
$ws.Range("H93").Value = 28867
$ws.Range("J93").Value = 28867
$ws.Range("L93").Value = 28867
$ws.Range("N93").Value = -33859
$ws.Range("H95").Value = 29966.666
$ws.Range("J95").Value = 29966.666
$ws.Range("L95").Value = 29966.666
$ws.Range("N95").Value = -35458.666
$ws.Range("H113").Value = 3080.0789
$ws.Range("I113").Value = 2408.4119
$ws.Range("J113").Value = 3623.8096
$ws.Range("K113").Value = 2408.4119
$ws.Range("L113").Value = 3623.8096
$ws.Range("M113").Value = 845.5880999999999
$ws.Range("N113").Value = -10131.8096

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8511.037
$ws.Range("I32").Value = 5357.9116
$ws.Range("J32").Value = 25004.309
$ws.Range("K32").Value = 5357.9116
$ws.Range("L32").Value = 25004.309
$ws.Range("M32").Value = -5070.9116
$ws.Range("N32").Value = -25578.309
$ws.Range("H61").Value = 2923.5833
$ws.Range("I61").Value = 1869
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 1869
$ws.Range("L61").Value = 4400
$ws.Range("M61").Value = -1657
$ws.Range("N61").Value = -4824
$ws.Range("H63").Value = 8088.909
$ws.Range("I63").Value = 8597.799999999999
$ws.Range("K63").Value = 8597.799999999999
$ws.Range("M63").Value = -7911.799999999999
$ws.Range("H66").Value = 8088.909
$ws.Range("I66").Value = 8597.799999999999
$ws.Range("K66").Value = 42989
$ws.Range("M66").Value = -39557
$ws.Range("H122").Value = 2745.6
$ws.Range("I122").Value = 2866
$ws.Range("J122").Value = 2063.3333
$ws.Range("K122").Value = 8598
$ws.Range("L122").Value = 6189.999899999999
$ws.Range("M122").Value = -6148
$ws.Range("N122").Value = -11089.9999
$ws.Range("H136").Value = 2923.5833
$ws.Range("I136").Value = 1869
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 5607
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -3057
$ws.Range("N136").Value = -18300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1213.8334
$ws.Range("I16").Value = 1382.8572
$ws.Range("J16").Value = 622.25
$ws.Range("K16").Value = 1382.8572
$ws.Range("L16").Value = 622.25
$ws.Range("M16").Value = -1095.8572
$ws.Range("N16").Value = -1196.25
$ws.Range("H58").Value = 2134.3235
$ws.Range("I58").Value = 1617
$ws.Range("J58").Value = 2716.3125
$ws.Range("K58").Value = 1617
$ws.Range("L58").Value = 2716.3125
$ws.Range("M58").Value = -1414
$ws.Range("N58").Value = -3122.3125
$ws.Range("H113").Value = 1213.8334
$ws.Range("I113").Value = 1382.8572
$ws.Range("J113").Value = 622.25
$ws.Range("K113").Value = 1382.8572
$ws.Range("L113").Value = 622.25
$ws.Range("M113").Value = 787.1428000000001
$ws.Range("N113").Value = -4962.25
$ws.Range("H122").Value = 1427.5
$ws.Range("I122").Value = 1660.8096
$ws.Range("J122").Value = 1139.2941
$ws.Range("K122").Value = 4982.4288
$ws.Range("L122").Value = 3417.8823
$ws.Range("M122").Value = -2532.4288
$ws.Range("N122").Value = -8317.882300000001
$ws.Range("H132").Value = 2165.2666
$ws.Range("I132").Value = 931.89655
$ws.Range("J132").Value = 4400.75
$ws.Range("K132").Value = 2795.68965
$ws.Range("L132").Value = 13202.25
$ws.Range("M132").Value = -265.6896500000003
$ws.Range("N132").Value = -18262.25
$ws.Range("H134").Value = 1489.1628
$ws.Range("I134").Value = 976.069
$ws.Range("J134").Value = 2552
$ws.Range("K134").Value = 2928.207
$ws.Range("L134").Value = 7656
$ws.Range("M134").Value = -393.2069999999999
$ws.Range("N134").Value = -12726
$ws.Range("H136").Value = 2134.3235
$ws.Range("I136").Value = 1617
$ws.Range("J136").Value = 2716.3125
$ws.Range("K136").Value = 4851
$ws.Range("L136").Value = 8148.9375
$ws.Range("M136").Value = -2301
$ws.Range("N136").Value = -13248.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 294.14285
$ws.Range("I50").Value = 212.2
$ws.Range("J50").Value = 499
$ws.Range("K50").Value = 636.5999999999999
$ws.Range("L50").Value = 1497
$ws.Range("M50").Value = -155.5999999999999
$ws.Range("N50").Value = -2459
$ws.Range("H53").Value = 294.14285
$ws.Range("I53").Value = 212.2
$ws.Range("J53").Value = 499
$ws.Range("K53").Value = 636.5999999999999
$ws.Range("L53").Value = 1497
$ws.Range("M53").Value = -155.5999999999999
$ws.Range("N53").Value = -2459
$ws.Range("H80").Value = 1055.8889
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1055.8889
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3167.6667
$ws.Range("N80").Value = -5039.6667
$ws.Range("H83").Value = 1055.8889
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1055.8889
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 9503.000099999999
$ws.Range("N83").Value = -18863.0001
$ws.Range("H97").Value = 350.34784
$ws.Range("I97").Value = 114.21429
$ws.Range("J97").Value = 717.6667
$ws.Range("K97").Value = 342.64287
$ws.Range("L97").Value = 2153.0001
$ws.Range("M97").Value = 153.35713
$ws.Range("N97").Value = -3145.0001
$ws.Range("H116").Value = 4543.3887
$ws.Range("I116").Value = 549.8
$ws.Range("K116").Value = 1649.4
$ws.Range("M116").Value = 1792.6
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3981.3333
$ws.Range("I80").Value = 4134
$ws.Range("J80").Value = 2760
$ws.Range("K80").Value = 4134
$ws.Range("L80").Value = 2760
$ws.Range("M80").Value = -3136
$ws.Range("N80").Value = -4756
$ws.Range("H83").Value = 3981.3333
$ws.Range("I83").Value = 4134
$ws.Range("J83").Value = 2760
$ws.Range("K83").Value = 20670
$ws.Range("L83").Value = 13800
$ws.Range("M83").Value = -15678
$ws.Range("N83").Value = -23784
$ws.Range("H122").Value = 1296.2307
$ws.Range("I122").Value = 1083.4445
$ws.Range("J122").Value = 1775
$ws.Range("K122").Value = 3250.3335
$ws.Range("L122").Value = 5325
$ws.Range("M122").Value = -800.3335000000002
$ws.Range("N122").Value = -10225

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14956.125
$ws.Range("I132").Value = 4300
$ws.Range("J132").Value = 21349.8
$ws.Range("K132").Value = 12900
$ws.Range("L132").Value = 64049.39999999999
$ws.Range("M132").Value = -10370
$ws.Range("N132").Value = -69109.39999999999
$ws.Range("H136").Value = 4578.6045
$ws.Range("I136").Value = 2449.2917
$ws.Range("K136").Value = 7347.875100000001
$ws.Range("M136").Value = -4797.875100000001

Write-Host "Applied all changes"